$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: extend table with two new columns (D, E) ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E14"))

# --- Step 2: header row - rename existing column, name the 2 new ones ---
$ws.Range("C1").Value = "Оценка контролно 1"
$ws.Range("D1").Value = "Оценка контролно 2"
$ws.Range("E1").Value = "Текуща оценка"

# --- Step 3: fill in faculty numbers for students who did not have one yet ---
$ws.Range("B7").Value = 1601651010
$ws.Range("B8").Value = 1601651003
$ws.Range("B11").Value = 1601651015
$ws.Range("B12").Value = 1601651013
$ws.Range("B13").Value = 1601651018

# --- Step 4: write new grade data (control 1 / control 2 / current grade) for every row, in original order ---
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 5

$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "N/A"
$ws.Range("D3").HorizontalAlignment = -4152
$ws.Range("E3").Value = 3

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 4.5

$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 4.5

$ws.Range("C6").Value = "N/A"
$ws.Range("C6").HorizontalAlignment = -4152
$ws.Range("D6").Value = "N/A"
$ws.Range("D6").HorizontalAlignment = -4152
$ws.Range("E6").Value = "N/A"
$ws.Range("E6").HorizontalAlignment = -4152

$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 4.5

$ws.Range("C8").Value = "N/A"
$ws.Range("C8").HorizontalAlignment = -4152
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 4

$ws.Range("C9").Value = "N/A"
$ws.Range("C9").HorizontalAlignment = -4152
$ws.Range("D9").Value = "N/A"
$ws.Range("D9").HorizontalAlignment = -4152
$ws.Range("E9").Value = "N/A"
$ws.Range("E9").HorizontalAlignment = -4152

$ws.Range("C10").Value = "N/A"
$ws.Range("C10").HorizontalAlignment = -4152
$ws.Range("D10").Value = "N/A"
$ws.Range("D10").HorizontalAlignment = -4152
$ws.Range("E10").Value = "N/A"
$ws.Range("E10").HorizontalAlignment = -4152

$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 4

$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 4.5

$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 4

# --- Step 5: sort the table by faculty number (col B), ascending, blanks last (matches Table1 sortCondition) ---
$ws.Range("A2:E13").Sort($ws.Range("B2"), 1, $null, $null, 1, $null, 1, 2)

# --- Step 6: set widths for the two new columns ---
$ws.Columns.Item(4).ColumnWidth = 25.42
$ws.Columns.Item(5).ColumnWidth = 16.6

# --- Step 7: reset selection to default top-left cell ---
$ws.Range("A1").Select()
